$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.747140288352966
$ws.Range("B1").Value = 2.634513854980469
$ws.Range("C1").Value = 3.333444595336914
$ws.Range("D1").Value = 2.174373865127563
$ws.Range("E1").Value = 0.5114942193031311
